# Insert a new weekly price record for "Coliflor" (Feria Lagunitas de Puerto
# Montt) as row 311, pushing the existing rows 311:332 down to 312:333
# (dimension grows from A1:R332 to A1:R333).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 311:332 down by inserting a fresh row at 311 (carries the D
# column's date-number-format style along, same as Excel's native
# "Insert Row" behaviour).
$ws.Rows.Item(311).Insert()

# Populate the newly inserted row with the new record.
$ws.Range("A311").Value = 4
$ws.Range("B311").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C311").Value = "Los Lagos"
$ws.Range("D311").Value = 44714
$ws.Range("E311").Value = 10
$ws.Range("F311").Value = 100112008
$ws.Range("G311").Value = "Coliflor"
$ws.Range("H311").Value = "Sin especificar"
$ws.Range("I311").Value = "Primera"
$ws.Range("J311").Value = 500
$ws.Range("K311").Value = 1500
$ws.Range("L311").Value = 1600
$ws.Range("M311").Value = 1550
$ws.Range("N311").Value = "$/unidad"
$ws.Range("O311").Value = "Región Metropolitana"
$ws.Range("P311").Value = 1550
$ws.Range("Q311").Value = 1
$ws.Range("R311").Value = "Hortaliza"
